$d = $word.ActiveDocument

function Split-Range($startPos, $endPos) {
    # Force run boundaries at both $startPos and $endPos by toggling Bold on
    # the (non-empty) range between them. Formatting is restored to its
    # original value, but the engine is left with the text on
    # [$startPos,$endPos) as its own <w:r>, distinct from its neighbours.
    $r = $d.Range($startPos, $endPos)
    $orig = $r.Bold
    $r.Bold = 1
    $r.Bold = $orig
}

# ---------------------------------------------------------------------------
# Change 1: ". We note that this manuscript is currently posted as a preprint
# on " -> ". " / "Please " / "note that this manuscript is currently posted
# as a preprint on "
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf(" We note that this manuscript is currently posted as a preprint on ")

# Replace "We " with "Please " (leaves the leading space just before it
# untouched); this necessarily merges the run with its neighbours, which we
# fix up below by re-splitting at the boundaries we actually want.
$weStart = $idx + 1
$weEnd = $idx + 4
$weRange = $d.Range($weStart, $weEnd)
$weRange.Text = "Please "

$p1a = $idx
$p1b = $idx + 1
$p1c = $idx + 8
Split-Range $p1a $p1b
Split-Range $p1b $p1c

# ---------------------------------------------------------------------------
# Change 2: "Here, we " -> "We" / " " (and "Here," is dropped)
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx2 = $full.IndexOf("Here, we ")

$hereStart = $idx2
$hereEnd = $idx2 + 9
$hereRange = $d.Range($hereStart, $hereEnd)
$hereRange.Text = "We "

$p2a = $idx2
$p2b = $idx2 + 2
Split-Range $p2a $p2b

# The text replacement above merges the (new) "We " run together with all of
# the following same-formatted runs in the paragraph ("explore whether ...",
# " ", "the", "ir statistical baselines", ". We build on ", "pioneer",
# "ing work by", " "). Re-impose the original run boundaries so only the
# "Here, we " -> "We"/" " run split shows up as a change.
$full = $d.Content.Text
$idx3 = $full.IndexOf("explore whether ecological SADs")
$offsets = @(0, 66, 67, 70, 94, 108, 115, 126, 127)
for ($i = 0; $i -lt $offsets.Length - 1; $i++) {
    $segStart = $idx3 + $offsets[$i]
    $segEnd = $idx3 + $offsets[$i + 1]
    Split-Range $segStart $segEnd
}
